# Fix to multiprocessing problem
#
# Renames the per-technology "SOURCE"/"EFF" header pair on each sheet to
# sheet-specific names (source_hs/eff_hs, source_cs/eff_cs, source_dhw/eff_dhw,
# source_el/eff_el), recodes the "SOURCE" categories into the new vocabulary
# (FUEL / GRID / SOLAR / PVT / DH / DC / PV), and - for COOLING and
# ELECTRICITY - inserts a new "eff" (COP / efficiency) column that used to be
# baked into the PEN/CO2 formulas as a literal divisor.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: structural changes first (insert the new efficiency columns) so
# that cross-sheet formulas written afterwards can target the final layout
# directly, and Excel's automatic reference-shifting on insert doesn't
# re-shift formulas we already wrote.
# ---------------------------------------------------------------------

# COOLING - insert a new "eff_cs" (COP) column between SOURCE and PEN
$wsCooling = $wb.Worksheets.Item("COOLING")
$wsCooling.Columns.Item(4).Insert()

# ELECTRICITY - insert new "source_el"/"eff_el" columns before PEN
$wsElectricity = $wb.Worksheets.Item("ELECTRICITY")
$wsElectricity.Range("C1:D1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# Step 2: DHW sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DHW")

$ws.Range("C1").Value = "source_dhw"
$ws.Range("D1").Value = "eff_dhw"

$ws.Range("C3").Value = "FUEL"
$ws.Range("C4").Value = "FUEL"
$ws.Range("C5").Value = "FUEL"
$ws.Range("C6").Value = "GRID"
$ws.Range("C7").Value = "FUEL"
$ws.Range("C8").Value = "GRID"
$ws.Range("C9").Value = "SOLAR"
$ws.Range("C10").Value = "GRID"
$ws.Range("C11").Value = "GRID"
$ws.Range("C12").Value = "DH"
$ws.Range("C13").Value = "DH"
$ws.Range("C14").Value = "DH"
$ws.Range("C15").Value = "DH"
$ws.Range("C16").Value = "PVT"
$ws.Range("C17").Value = "DH"
$ws.Range("C18").Value = "DH"
$ws.Range("C19").Value = "DH"
$ws.Range("C20").Value = "DH"
$ws.Range("C21").Value = "DH"
$ws.Range("C22").Value = "DH"

$ws.Range("E10").Formula = "=ELECTRICITY!E3*1.15/2.7"
$ws.Range("F10").Formula = "=1.15*ELECTRICITY!F3/2.7"

$ws.Range("E22").Formula = "=(0.43*(E21)+0.28*(ELECTRICITY!E3/2.96)+0.18*0.954+0.11*0)"
$ws.Range("F22").Formula = "=(0.43*(F21)+0.28*(ELECTRICITY!F3/2.96)+0.18*0.0149+0.11*0)"

$ws.Columns.Item(4).ColumnWidth = 8.14
$ws.Range("C9").Select()

# ---------------------------------------------------------------------
# Step 3: HEATING sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("HEATING")

$ws.Range("C1").Value = "source_hs"
$ws.Range("D1").Value = "eff_hs"

$ws.Range("C3").Value = "FUEL"
$ws.Range("C4").Value = "FUEL"
$ws.Range("C5").Value = "FUEL"
$ws.Range("C6").Value = "GRID"
$ws.Range("C7").Value = "FUEL"
$ws.Range("C8").Value = "GRID"
$ws.Range("C9").Value = "SOLAR"
$ws.Range("C10").Value = "GRID"
$ws.Range("C11").Value = "GRID"
$ws.Range("C12").Value = "DH"
$ws.Range("C13").Value = "DH"
$ws.Range("C14").Value = "DH"
$ws.Range("C15").Value = "DH"
$ws.Range("C16").Value = "PVT"
$ws.Range("C17").Value = "DH"
$ws.Range("C18").Value = "DH"
$ws.Range("C19").Value = "DH"
$ws.Range("C20").Value = "DH"
$ws.Range("C21").Value = "DH"
$ws.Range("C22").Value = "DH"

$ws.Range("E10").Formula = "=ELECTRICITY!E3*1.15/2.7"
$ws.Range("F10").Formula = "=ELECTRICITY!F3*1.15/2.7"

$ws.Range("E22").Formula = "=(0.43*(E21)+0.28*(ELECTRICITY!E3/2.96)+0.18*0.954+0.11*0)"
$ws.Range("F22").Formula = "=(0.43*(F21)+0.28*(ELECTRICITY!F3/2.96)+0.18*0.0149+0.11*0)"

$ws.Range("F7").Select()

# ---------------------------------------------------------------------
# Step 4: COOLING sheet - fill in the new column + update source codes
# ---------------------------------------------------------------------
$ws = $wsCooling

$ws.Range("C1").Value = "source_cs"
$ws.Range("D1").Value = "eff_cs"

$ws.Range("C2").Value = "none"
$ws.Range("D2").Value = 0

$ws.Range("C3").Value = "GRID"
$ws.Range("D3").Value = 2.7

$ws.Range("C4").Value = "GRID"
$ws.Range("D4").Value = 3

$ws.Range("C5").Value = "DC"
$ws.Range("D5").Value = 3.2

$ws.Range("C6").Value = "DC"
$ws.Range("D6").Value = 2.8

$ws.Range("E6").Formula = "=ELECTRICITY!E8/4"
$ws.Range("F6").Formula = "=ELECTRICITY!F8/4"
$ws.Range("G6").Formula = "=ELECTRICITY!G8/4"

$ws.Columns.Item(4).ColumnWidth = 9.14
$ws.Range("C5").Select()

# ---------------------------------------------------------------------
# Step 5: ELECTRICITY sheet - fill in the new source/eff columns
# ---------------------------------------------------------------------
$ws = $wsElectricity

$ws.Range("C1").Value = "source_el"
$ws.Range("D1").Value = "eff_el"

$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

$ws.Range("C3").Value = "GRID"
$ws.Range("D3").Value = 0.99

$ws.Range("C4").Value = "PV"
$ws.Range("D4").Value = 0.99

$ws.Range("C5").Value = "PVT"
$ws.Range("D5").Value = 0.99

$ws.Range("C6").Value = "GRID"
$ws.Range("D6").Value = 0.99

$ws.Range("C7").Value = "GRID"
$ws.Range("D7").Value = 0.99

$ws.Range("C8").Value = "GRID"
$ws.Range("D8").Value = 0.99

$ws.Columns.Item(3).ColumnWidth = 9.14
$ws.Columns.Item(4).ColumnWidth = 9.14
$ws.Range("F16").Select()
